$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.659095525741577
$ws.Range("B1").Value = 2.530667781829834
$ws.Range("C1").Value = 2.925381660461426
$ws.Range("D1").Value = 3.173591136932373
$ws.Range("E1").Value = 1.054108381271362
